# "Add files via upload" — re-upload of SEO Project.pptx with PowerPoint's
# own auto-refresh touch-ups baked in:
#   1. Every "Date Placeholder" footer field (11 slide layouts + the slide
#      master) gets its cached datetimeFigureOut text bumped from
#      23-08-2024 to 24-08-2024 (the file was re-saved a day later).
#   2. The leftover empty "Content Placeholder 2" autolayout box on the
#      final "Thank You" slide is removed.

$p = $ppt.ActivePresentation

$oldDate = "23-08-2024"
$newDate = "24-08-2024"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1a. Slide master's Date Placeholder.
$master = $p.SlideMaster
Update-DateShape $master.Shapes

# 1b. Every slide layout's Date Placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShape $layout.Shapes
}

# 2. Drop the empty "Content Placeholder 2" box left on the last slide.
$lastSlide = $p.Slides.Item($p.Slides.Count)
for ($i = $lastSlide.Shapes.Count; $i -ge 1; $i--) {
    $shp = $lastSlide.Shapes.Item($i)
    if ($shp.Name -eq "Content Placeholder 2") {
        $shp.Delete()
    }
}
